# Refactor with sew to implement code chunks:
# Row 1 now holds the R source of the code chunk (rendered in a
# monospace font), a blank spacer row follows, and the original
# printed output (the named vector a/b/c -> 1/2/3) is shifted down
# two rows to rows 3-5.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Make room: push the existing output down by two rows (row 1 -> row 3,
# row 2 -> row 4, row 3 -> row 5), leaving a blank row 2 as a spacer.
$ws.Rows("1:2").Insert()

# Code-chunk text goes in the newly freed A1, rendered in a monospace font.
$ws.Range("A1").Value = "setNames(1:3, c('a', 'b', 'c'))"
$ws.Range("A1").Font.Name = "Courier New"
